# Weekly data refresh: insert a new observation at row 59 (pushing the
# existing 59-154 block down to 60-155, growing the used range from
# A1:R154 to A1:R155) and populate the new row with the latest reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the current row 59; everything below shifts
# down by one (old 59 -> 60, old 60 -> 61, ..., old 154 -> 155).
$ws.Rows.Item(59).Insert()

# Populate the newly inserted row 59 with the new record.
$ws.Range("A59").Value = 10
$ws.Range("B59").Value = "Vega Modelo de Temuco"
$ws.Range("C59").Value = "La Araucanía"
$ws.Range("D59").Value = 44665
$ws.Range("E59").Value = 9
$ws.Range("F59").Value = 100114007
$ws.Range("G59").Value = "Jengibre"
$ws.Range("H59").Value = "Sin especificar"
$ws.Range("I59").Value = "Primera"
$ws.Range("J59").Value = 60
$ws.Range("K59").Value = 15000
$ws.Range("L59").Value = 25000
$ws.Range("M59").Value = 21667
$ws.Range("N59").Value = "$/caja 13 kilos"
$ws.Range("O59").Value = "Perú"
$ws.Range("P59").Value = 1667
$ws.Range("Q59").Value = 13
$ws.Range("R59").Value = "Hortaliza"
